# Adding "Area"/"Atotal" columns (G/H) to the discharge worksheet, mirroring
# the existing "Q"/"Qtotal" (E/F) columns but using (depth*width) instead of
# (depth*width*velocity).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells for the second table (row 11) ------------------------
$ws.Range("G11").Value = "Area"
$ws.Range("H11").Value = "Atotal"

# --- New "Area" formulas down column G --------------------------------------
$ws.Range("G12").Formula = "=(D12-0)*B12/100"
$ws.Range("G13").Formula = "=(D13-D12)*B13/100"
$ws.Range("G14:G25").Formula = "=(D14-D13)*B14/100"

# --- New "Atotal" formula in H12 --------------------------------------------
$ws.Range("H12").Formula = "=SUM(G12:G21)"

# --- Best-effort re-alignment of existing direct formatting ----------------
# The workbook's style table originally carried a redundant duplicate xf
# (fontId 1 / fillId 0, once with a no-op applyFill="1" and once without).
# Re-applying the named cell style lets the engine intern cells onto the
# single canonical xf instead of the redundant one. Only previously-populated
# cells are touched (one area at a time - a multi-area Range only restyles
# its first area), so no blank cells get materialised along the way.
$restyleAddrs = @(
    "A1:E1", "A2:E2",
    "A3", "C3:E3",
    "A4", "C4:E4",
    "A5", "C5:E5",
    "A6", "C6:E6",
    "A7", "C7:E7",
    "A8", "C8:E8",
    "A12:C12",
    "A13:C18"
)
foreach ($addr in $restyleAddrs) {
    $ws.Range($addr).Style = "Normal 2"
}

# --- View state: selection moved to H12, top-left cell scrolled to B5 ------
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("H12").Select()
